# DataImport.xlsx — "Add files via upload" re-edit
# Net effect of the commit (once the shared-string table index churn is
# resolved): three data cells got new text, the data block A1:F3 gained a
# thin box border, the saved selection moved to M13, and the sheet picked
# up an explicit (default) page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content changes -------------------------------------------------
# B2: 23-3  -> 23-33
# B3: 23-3  -> 23-34
# E3: 23-123148 -> 23-12314
$ws.Range("B2").Value = "23-33"
$ws.Range("B3").Value = "23-34"
$ws.Range("E3").Value = "23-12314"

# --- Formatting: thin box border around the whole data range -------------
$ws.Range("A1:F3").Borders.LineStyle = 1   # xlContinuous
$ws.Range("A1:F3").Borders.Weight = 2      # xlThin

# --- Selection moved to M13 (saved cursor position) -----------------------
$ws.Range("M13").Select()

# --- Page setup (paper size / orientation became explicit) ---------------
$ws.PageSetup.PaperSize = 9    # xlPaperA4
$ws.PageSetup.Orientation = 1  # xlPortrait
